$d = $word.ActiveDocument

$pairs = @(
    @{old="263×6=1578"; new="257×2=514"},
    @{old="203×4=812"; new="154×5=770"},
    @{old="453×4=1812"; new="978×5=4890"},
    @{old="549×2=1098"; new="844×9=7596"},
    @{old="486×6=2916"; new="966×2=1932"},
    @{old="977×8=7816"; new="170×7=1190"},
    @{old="759×8=6072"; new="854×7=5978"},
    @{old="855×7=5985"; new="259×4=1036"},
    @{old="729×9=6561"; new="614×3=1842"},
    @{old="839×2=1678"; new="167×3=501"},
    @{old="354×2=708"; new="760×7=5320"},
    @{old="483×4=1932"; new="380×9=3420"},
    @{old="350×8=2800"; new="936×4=3744"},
    @{old="525×4=2100"; new="820×3=2460"},
    @{old="455×6=2730"; new="390×2=780"},
    @{old="743×4=2972"; new="236×8=1888"},
    @{old="882×4=3528"; new="533×3=1599"},
    @{old="635×8=5080"; new="724×2=1448"},
    @{old="292×4=1168"; new="181×7=1267"},
    @{old="859×3=2577"; new="924×9=8316"},
    @{old="845×9=7605"; new="797×4=3188"},
    @{old="791×7=5537"; new="128×7=896"},
    @{old="311×3=933"; new="865×9=7785"},
    @{old="809×7=5663"; new="106×4=424"},
    @{old="903×3=2709"; new="807×6=4842"}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
